# Apply crypto price/volume updates to match the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.299.73"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.830.99"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'626.93"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").Value = "'166.33"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "3.829.16"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'6.59"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "'36.20"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "4.474.54"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.283.31"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.594.86"
$ws.Range("E17").Value = "  -5.82%  "
$ws.Range("D18").Value = "'18.13"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "'7.16"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'467.74"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'9.70"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'0.710"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "'0.0000154"
$ws.Range("E24").Value = "  +4.81%  "
$ws.Range("D25").Value = "'83.94"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'12.05"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'2.16"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "3.976.79"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D34").Value = "'29.31"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'9.13"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  +7.54%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +5.52%  "
$ws.Range("D40").Value = "'5.94"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("D45").Value = "'0.300"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "'154.83"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "'46.99"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'42.65"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.48"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").Value = "'0.000278"
